$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.572.58'
$ws.Range("E2").Value = '  +4.74%  '
$ws.Range("D3").Value = '3.632.68'
$ws.Range("E3").Value = '  +4.34%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '201.43'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +8.97%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '583.96'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.78%  '
$ws.Range("D7").Value = '3.627.41'
$ws.Range("E7").Value = '  +4.27%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.624'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.22%  '
$ws.Range("E9").Value = '  +0.08%  '
$ws.Range("E10").Value = '  +6.18%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '60.80'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +18.02%  '
$ws.Range("E12").Value = '  +6.80%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000287'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +14.50%  '
$ws.Range("E14").Value = '  +6.70%  '
$ws.Range("D15").Value = '4.206.06'
$ws.Range("E15").Value = '  +4.03%  '
$ws.Range("D16").Value = '3.638.86'
$ws.Range("E16").Value = '  +4.28%  '
$ws.Range("E17").Value = '  +1.35%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '19.26'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +7.91%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.59'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.80%  '
$ws.Range("D20").Value = '68.486.07'
$ws.Range("E20").Value = '  +4.94%  '
$ws.Range("E21").Value = '  +4.68%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '407.16'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.59%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.05'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +21.56%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.32'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.71%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '86.00'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.19%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.94'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.57%  '
$ws.Range("E27").Value = '  +17.33%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.77'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.45%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.15'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.26%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.49'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +10.26%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.85'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +12.84%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '31.86'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.39%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '683.19'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +12.29%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '12.33'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.17%  '
$ws.Range("E35").Value = '  +4.62%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '63.93'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.35%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '42.13'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.35%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.420'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.69%  '
$ws.Range("E39").Value = '  +0.10%  '
$ws.Range("D40").Value = '0.0₃0777'
$ws.Range("E40").Value = '  +6.43%  '
$ws.Range("E41").Value = '  +19.01%  '
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.136'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.40%  '
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '3.196.02'
$ws.Range("E43").Value = '  +8.89%  '
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.997'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.28%  '
$ws.Range("B45").Value = 'Fetch.AI'
$ws.Range("C45").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.71'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +12.64%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.87'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +27.61%  '
$ws.Range("E47").Value = '  +16.80%  '
$ws.Range("E48").Value = '  +6.89%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.132'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.96%  '
$ws.Range("E50").Value = '  +8.04%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.09'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.94%  '
